$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Text number format to the Price (D) column data rows so that
# numeric-looking price strings (e.g. "3.47") are stored as literal text
# instead of being coerced into floating point numbers, matching the
# original inlineStr cell type. ClearFormats() afterwards removes the
# temporary style again so cells keep their original (default) styling.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range('D2').Value = '94.971.99'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '3.442.61'
$ws.Range('E3').Value = '  +3.37%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '239.20'
$ws.Range('E5').Value = '  -4.33%  '
$ws.Range('D6').Value = '641.33'
$ws.Range('E6').Value = '  -2.20%  '
$ws.Range('D7').Value = '1.46'
$ws.Range('E7').Value = '  +4.65%  '
$ws.Range('D8').Value = '0.399'
$ws.Range('E8').Value = '  -5.08%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = '0.991'
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('D11').Value = '3.438.88'
$ws.Range('E11').Value = '  +3.33%  '
$ws.Range('D12').Value = '0.197'
$ws.Range('E12').Value = '  -4.30%  '
$ws.Range('D13').Value = '41.28'
$ws.Range('E13').Value = '  +2.51%  '
$ws.Range('D14').Value = '6.06'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = '94.850.37'
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('D16').Value = '4.093.72'
$ws.Range('E16').Value = '  +3.48%  '
$ws.Range('D17').Value = '0.0000255'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '8.42'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '3.435.79'
$ws.Range('E19').Value = '  +2.81%  '
$ws.Range('D20').Value = '17.69'
$ws.Range('E20').Value = '  +3.51%  '
$ws.Range('D21').Value = '11.42'
$ws.Range('E21').Value = '  +8.36%  '
$ws.Range('D22').Value = '0.506'
$ws.Range('E22').Value = '  -5.82%  '
$ws.Range('D23').Value = '500.57'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('D24').Value = '3.15'
$ws.Range('E24').Value = '  -5.83%  '
$ws.Range('D25').Value = '0.0000190'
$ws.Range('E25').Value = '  -3.60%  '
$ws.Range('D26').Value = '6.54'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').Value = '91.79'
$ws.Range('E27').Value = '  -4.74%  '
$ws.Range('D28').Value = '3.632.01'
$ws.Range('E28').Value = '  +3.54%  '
$ws.Range('D29').Value = '11.94'
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('D30').Value = '11.66'
$ws.Range('E30').Value = '  +5.38%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = '2.74'
$ws.Range('E32').Value = '  +10.38%  '
$ws.Range('D33').Value = '0.136'
$ws.Range('E33').Value = '  -4.84%  '
$ws.Range('D34').Value = '0.183'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('D35').Value = '30.98'
$ws.Range('E35').Value = '  +10.74%  '
$ws.Range('E36').Value = '  +0.48%  '
$ws.Range('D37').Value = '0.564'
$ws.Range('E37').Value = '  +3.36%  '
$ws.Range('D38').Value = '7.65'
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('D39').Value = '1.43'
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('D40').Value = '522.76'
$ws.Range('E40').Value = '  +2.84%  '
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').Value = '0.150'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('D43').Value = '0.907'
$ws.Range('E43').Value = '  +9.07%  '
$ws.Range('D44').Value = '24.09'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('D45').Value = '1.69'
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('D46').Value = '5.58'
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('D47').Value = '0.0414'
$ws.Range('E47').Value = '  -2.58%  '
$ws.Range('B48').Value = 'MantraDAO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D48').Value = '3.47'
$ws.Range('E48').Value = '  -4.93%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '53.50'
$ws.Range('E49').Value = '  +0.42%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '2.14'
$ws.Range('E50').Value = '  +8.21%  '
$ws.Range('D51').Value = '3.16'
$ws.Range('E51').Value = '  +1.26%  '

$priceCol.ClearFormats()
